$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Daily_Data")

# New daily report block: 2026-02-19 (serial 46072), published 2026-02-20.
# Column layout: A=Date, B=Region_Type, C=PREV_TOTAL, D=RECEIVED, E=WITHDRAWN,
# F=NET_CHANGE, G=ADJUSTMENT, H=TOTAL_TODAY.
$reportDate = 46072
$startRow = 706

$rows = @(
  @('ASAHI DEPOSITORY LLC Registered', 0, 0, 0, 0, 0, 0),
  @('ASAHI DEPOSITORY LLC Eligible', 0, 0, 0, 0, 0, 0),
  @('BRINK''S, INC. Registered', 71275.599, 0, 0, 0, 0, 71275.599),
  @('BRINK''S, INC. Eligible', 86539.92200000001, 0, 0, 0, 0, 86539.92200000001),
  @('CNT DEPOSITORY, INC. Registered', 1246.06, 0, 0, 0, 0, 1246.06),
  @('CNT DEPOSITORY, INC. Eligible', 0, 0, 0, 0, 0, 0),
  @('DELAWARE DEPOSITORY Registered', 1633.941, 0, 0, 0, 0, 1633.941),
  @('DELAWARE DEPOSITORY Eligible', 18459.584, 0, 0, 0, 0, 18459.584),
  @('HSBC BANK, USA Registered', 1394.758, 0, 0, 0, 0, 1394.758),
  @('HSBC BANK, USA Eligible', 9281.978999999999, 0, 0, 0, 0, 9281.978999999999),
  @('INTERNATIONAL DEPOSITORY SERVICES OF DELAWARE Registered', 2395.448, 0, 0, 0, 0, 2395.448),
  @('INTERNATIONAL DEPOSITORY SERVICES OF DELAWARE Eligible', 0, 0, 0, 0, 0, 0),
  @('JP MORGAN CHASE BANK NA Registered', 113601.788, 0, 0, 0, 0, 113601.788),
  @('JP MORGAN CHASE BANK NA Eligible', 75944.144, 0, 0, 0, 0, 75944.144),
  @('LOOMIS INTERNATIONAL (US) LLC Registered', 59209.788, 0, 0, 0, 0, 59209.788),
  @('LOOMIS INTERNATIONAL (US) LLC Eligible', 70953.296, 0, 0, 0, 0, 70953.296),
  @('MALCA-AMIT USA, LLC Registered', 395.145, 0, 0, 0, 0, 395.145),
  @('MALCA-AMIT USA, LLC Eligible', 0, 0, 0, 0, 0, 0),
  @('MANFRA, TORDELLA & BROOKES, LLC Registered', 48292.647, 0, 0, 0, 0, 48292.647),
  @('MANFRA, TORDELLA & BROOKES, LLC Eligible', 3432.284, 0, 0, 0, 0, 3432.284),
  @('STONEX PRECIOUS METALS LLC Registered', 14122.765, 0, 0, 0, 0, 14122.765),
  @('STONEX PRECIOUS METALS LLC Eligible', 16.075, 0, 0, 0, 0, 16.075)
)

$r = $startRow
foreach ($row in $rows) {
  $ws.Cells.Item($r, 1).Value = $reportDate
  $ws.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
  $ws.Cells.Item($r, 2).Value = $row[0]
  $ws.Cells.Item($r, 3).Value = $row[1]
  $ws.Cells.Item($r, 4).Value = $row[2]
  $ws.Cells.Item($r, 5).Value = $row[3]
  $ws.Cells.Item($r, 6).Value = $row[4]
  $ws.Cells.Item($r, 7).Value = $row[5]
  $ws.Cells.Item($r, 8).Value = $row[6]
  $r = $r + 1
}
